# Update "想去人数" (want-to-go count) figures in both the "展览" sheet
# and the consolidated "全部类型" sheet, per the site's latest scrape.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 5185
$wsExpo.Range("F5").Value = 7489
$wsExpo.Range("F12").Value = 4334
$wsExpo.Range("F16").Value = 2932
$wsExpo.Range("F18").Value = 568
$wsExpo.Range("F19").Value = 212
$wsExpo.Range("F20").Value = 510
$wsExpo.Range("F21").Value = 446
$wsExpo.Range("F22").Value = 464
$wsExpo.Range("F23").Value = 313
$wsExpo.Range("F28").Value = 1390
$wsExpo.Range("F29").Value = 110
$wsExpo.Range("F30").Value = 584
$wsExpo.Range("F37").Value = 2936
$wsExpo.Range("F38").Value = 710
$wsExpo.Range("F40").Value = 91
$wsExpo.Range("F42").Value = 47

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5185
$wsAll.Range("F5").Value = 7489
$wsAll.Range("F12").Value = 4334
$wsAll.Range("F16").Value = 2932
$wsAll.Range("F18").Value = 568
$wsAll.Range("F19").Value = 212
$wsAll.Range("F20").Value = 510
$wsAll.Range("F21").Value = 446
$wsAll.Range("F22").Value = 464
$wsAll.Range("F24").Value = 313
$wsAll.Range("F29").Value = 1390
$wsAll.Range("F30").Value = 110
$wsAll.Range("F31").Value = 584
$wsAll.Range("F38").Value = 2936
$wsAll.Range("F40").Value = 710
$wsAll.Range("F42").Value = 91
$wsAll.Range("F44").Value = 47
